$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 823797.3
$ws.Range("J17").Value = 849531.6
$ws.Range("L17").Value = 2548594.8
$ws.Range("N17").Value = -2548930.8
$ws.Range("H68").Value = 17999.8
$ws.Range("J68").Value = 17999.8
$ws.Range("L68").Value = 17999.8
$ws.Range("N68").Value = -19497.8
$ws.Range("H71").Value = 17999.8
$ws.Range("J71").Value = 17999.8
$ws.Range("L71").Value = 53999.39999999999
$ws.Range("N71").Value = -61487.39999999999
$ws.Range("H82").Value = 4098.3125
$ws.Range("I82").Value = 2577.4285
$ws.Range("J82").Value = 5281.222
$ws.Range("K82").Value = 7732.2855
$ws.Range("L82").Value = 15843.666
$ws.Range("M82").Value = -7326.2855
$ws.Range("N82").Value = -16655.666
$ws.Range("H85").Value = 4098.3125
$ws.Range("I85").Value = 2577.4285
$ws.Range("J85").Value = 5281.222
$ws.Range("K85").Value = 7732.2855
$ws.Range("L85").Value = 15843.666
$ws.Range("M85").Value = -6328.2855
$ws.Range("N85").Value = -18651.666
$ws.Range("H93").Value = 36993.8
$ws.Range("J93").Value = 36993.8
$ws.Range("L93").Value = 36993.8
$ws.Range("N93").Value = -41985.8
$ws.Range("H127").Value = 415172.8
$ws.Range("I127").Value = 343.57144
$ws.Range("J127").Value = 568004.6
$ws.Range("K127").Value = 1030.71432
$ws.Range("L127").Value = 1704013.8
$ws.Range("M127").Value = 3929.28568
$ws.Range("N127").Value = -1713933.8
$ws.Range("H138").Value = 12823656
$ws.Range("I138").Value = 22224956
$ws.Range("J138").Value = 3703.0908
$ws.Range("K138").Value = 66674868
$ws.Range("L138").Value = 11109.2724
$ws.Range("M138").Value = -66669728
$ws.Range("N138").Value = -21389.2724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2067324.2
$ws.Range("I45").Value = 3031475.5
$ws.Range("J45").Value = 1285.7142
$ws.Range("K45").Value = 3031475.5
$ws.Range("L45").Value = 1285.7142
$ws.Range("M45").Value = -3031098.5
$ws.Range("N45").Value = -2039.7142
$ws.Range("H61").Value = 19233664
$ws.Range("I61").Value = 21742028
$ws.Range("J61").Value = 2872.6667
$ws.Range("K61").Value = 21742028
$ws.Range("L61").Value = 2872.6667
$ws.Range("M61").Value = -21741816
$ws.Range("N61").Value = -3296.6667
$ws.Range("H102").Value = 2859.6
$ws.Range("I102").Value = 3000
$ws.Range("K102").Value = 3000
$ws.Range("M102").Value = -1378
$ws.Range("H117").Value = 23907
$ws.Range("J117").Value = 23907
$ws.Range("L117").Value = 23907
$ws.Range("N117").Value = -33085
$ws.Range("H136").Value = 19233664
$ws.Range("I136").Value = 21742028
$ws.Range("J136").Value = 2872.6667
$ws.Range("K136").Value = 65226084
$ws.Range("L136").Value = 8618.000100000001
$ws.Range("M136").Value = -65223534
$ws.Range("N136").Value = -13718.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 31252856
$ws.Range("I86").Value = 2719.9
$ws.Range("J86").Value = 83336420
$ws.Range("K86").Value = 2719.9
$ws.Range("L86").Value = 83336420
$ws.Range("M86").Value = -1596.9
$ws.Range("N86").Value = -83338666
$ws.Range("H89").Value = 31252856
$ws.Range("I89").Value = 2719.9
$ws.Range("J89").Value = 83336420
$ws.Range("K89").Value = 13599.5
$ws.Range("L89").Value = 416682100
$ws.Range("M89").Value = -7983.5
$ws.Range("N89").Value = -416693332
$ws.Range("H105").Value = 3162.6223
$ws.Range("I105").Value = 1637.2727
$ws.Range("J105").Value = 4621.6523
$ws.Range("K105").Value = 1637.2727
$ws.Range("L105").Value = 4621.6523
$ws.Range("M105").Value = 109.7273
$ws.Range("N105").Value = -8115.6523
$ws.Range("H118").Value = 7684.375
$ws.Range("J118").Value = 7684.375
$ws.Range("L118").Value = 7684.375
$ws.Range("N118").Value = -10998.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7580121
$ws.Range("I31").Value = 4710.8438
$ws.Range("J31").Value = 27781216
$ws.Range("K31").Value = 4710.8438
$ws.Range("L31").Value = 27781216
$ws.Range("M31").Value = -4415.8438
$ws.Range("N31").Value = -27781806
$ws.Range("H34").Value = 7580121
$ws.Range("I34").Value = 4710.8438
$ws.Range("J34").Value = 27781216
$ws.Range("K34").Value = 4710.8438
$ws.Range("L34").Value = 27781216
$ws.Range("M34").Value = -4508.8438
$ws.Range("N34").Value = -27781620
$ws.Range("H62").Value = 2916.6667
$ws.Range("I62").Value = 2200
$ws.Range("J62").Value = 3275
$ws.Range("K62").Value = 2200
$ws.Range("L62").Value = 3275
$ws.Range("M62").Value = -1576
$ws.Range("N62").Value = -4523
$ws.Range("H65").Value = 2916.6667
$ws.Range("I65").Value = 2200
$ws.Range("J65").Value = 3275
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 16375
$ws.Range("M65").Value = -7880
$ws.Range("N65").Value = -22615
$ws.Range("H99").Value = 1759.75
$ws.Range("I99").Value = 1686.3334
$ws.Range("J99").Value = 1980
$ws.Range("K99").Value = 1686.3334
$ws.Range("L99").Value = 1980
$ws.Range("M99").Value = -188.3334
$ws.Range("N99").Value = -4976
$ws.Range("H126").Value = 1759.75
$ws.Range("I126").Value = 1686.3334
$ws.Range("J126").Value = 1980
$ws.Range("K126").Value = 5059.0002
$ws.Range("L126").Value = 5940
$ws.Range("M126").Value = -2589.0002
$ws.Range("N126").Value = -10880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 85714350
$ws.Range("I33").Value = 100000040
$ws.Range("K33").Value = 600000240
$ws.Range("M33").Value = -599999957
$ws.Range("H109").Value = 4085.4167
$ws.Range("I109").Value = 2000
$ws.Range("J109").Value = 4383.3335
$ws.Range("K109").Value = 6000
$ws.Range("L109").Value = 13150.0005
$ws.Range("M109").Value = -4960
$ws.Range("N109").Value = -15230.0005
$ws.Range("H116").Value = 615.6
$ws.Range("I116").Value = 519.5
$ws.Range("K116").Value = 1558.5
$ws.Range("M116").Value = 1883.5
$ws.Range("H121").Value = 568.4286
$ws.Range("I121").Value = 353.33334
$ws.Range("J121").Value = 855.2222
$ws.Range("K121").Value = 1060.00002
$ws.Range("L121").Value = 2565.6666
$ws.Range("M121").Value = 249.9999800000001
$ws.Range("N121").Value = -5185.6666
$ws.Range("H131").Value = 1465.6578
$ws.Range("J131").Value = 1170.5358
$ws.Range("L131").Value = 3511.6074
$ws.Range("N131").Value = -13591.6074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12094.435
$ws.Range("I70").Value = 20911.273
$ws.Range("J70").Value = 4012.3333
$ws.Range("K70").Value = 20911.273
$ws.Range("L70").Value = 4012.3333
$ws.Range("M70").Value = -20641.273
$ws.Range("N70").Value = -4552.3333
$ws.Range("H73").Value = 12094.435
$ws.Range("I73").Value = 20911.273
$ws.Range("J73").Value = 4012.3333
$ws.Range("K73").Value = 20911.273
$ws.Range("L73").Value = 4012.3333
$ws.Range("M73").Value = -19975.273
$ws.Range("N73").Value = -5884.3333
$ws.Range("H118").Value = 14600
$ws.Range("J118").Value = 14600
$ws.Range("L118").Value = 14600
$ws.Range("N118").Value = -17914

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3899.4075
$ws.Range("I40").Value = 5515.909
$ws.Range("J40").Value = 2788.0625
$ws.Range("K40").Value = 5515.909
$ws.Range("L40").Value = 2788.0625
$ws.Range("M40").Value = -5379.909
$ws.Range("N40").Value = -3060.0625
$ws.Range("H55").Value = 173.1875
$ws.Range("I55").Value = 78.25
$ws.Range("J55").Value = 204.83333
$ws.Range("K55").Value = 78.25
$ws.Range("L55").Value = 204.83333
$ws.Range("M55").Value = 94.75
$ws.Range("N55").Value = -550.8333299999999
$ws.Range("H61").Value = 1525
$ws.Range("I61").Value = 1033.3334
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1033.3334
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -831.3334
$ws.Range("N61").Value = -3404
$ws.Range("H100").Value = 2551.375
$ws.Range("I100").Value = 2295.6667
$ws.Range("J100").Value = 2704.8
$ws.Range("K100").Value = 2295.6667
$ws.Range("L100").Value = 2704.8
$ws.Range("M100").Value = -1754.6667
$ws.Range("N100").Value = -3786.8
$ws.Range("H113").Value = 1525
$ws.Range("I113").Value = 1033.3334
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1033.3334
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1136.6666
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 11367882
$ws.Range("I132").Value = 3750.8462
$ws.Range("J132").Value = 27782740
$ws.Range("K132").Value = 11252.5386
$ws.Range("L132").Value = 83348220
$ws.Range("M132").Value = -8722.5386
$ws.Range("N132").Value = -83353280
